$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.165518666666667
$row[0,7] = 3.496556
$row[0,8] = 0.1888647065994748
$row[0,9] = 0.1888647065994748
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.063353333333333
$row[0,13] = 9.190059999999999
$row[0,14] = 0.1884019917097105
$row[0,15] = 0.1884019917097105
$row[0,16] = 3.570395492595555
$row[0,17] = 32.13355943336
$row[0,18] = 0.03558248688701115
$row[0,19] = 0.03558248688701115
$ws.Range("A2:T2").Value = $row

# Row 3
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.165518666666667
$row[0,7] = 3.496556
$row[0,8] = 0.1888647065994748
$row[0,9] = 0.1888647065994748
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 6.757657999999999
$row[0,13] = 20.272974
$row[0,14] = 0.4156086771445645
$row[0,15] = 0.4156086771445645
$row[0,16] = 7.876176541949332
$row[0,17] = 70.885588877544
$row[0,18] = 0.078493810869104
$row[0,19] = 0.07849381086910401
$ws.Range("A3:T3").Value = $row

# Row 4
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "MuSCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.165518666666667
$row[0,7] = 3.496556
$row[0,8] = 0.1888647065994748
$row[0,9] = 0.1888647065994748
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.493414666666666
$row[0,13] = 10.480244
$row[0,14] = 0.214851572590793
$row[0,15] = 0.214851572590793
$row[0,16] = 4.071640004407111
$row[0,17] = 36.644760039664
$row[0,18] = 0.04057787921979587
$row[0,19] = 0.04057787921979588
$ws.Range("A4:T4").Value = $row

# Row 5
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "Resolving-Mac"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.165518666666667
$row[0,7] = 3.496556
$row[0,8] = 0.1888647065994748
$row[0,9] = 0.1888647065994748
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.945239333333333
$row[0,13] = 8.835718
$row[0,14] = 0.1811377585549322
$row[0,15] = 0.1811377585549322
$row[0,16] = 3.432731420800889
$row[0,17] = 30.894582787208
$row[0,18] = 0.03421052962356377
$row[0,19] = 0.03421052962356377
$ws.Range("A5:T5").Value = $row

# Row 6
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "ECs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 0.3739756666666667
$row[0,7] = 1.121927
$row[0,8] = 0.06060032033836409
$row[0,9] = 0.06060032033836409
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.063353333333333
$row[0,13] = 9.190059999999999
$row[0,14] = 0.1884019917097105
$row[0,15] = 0.1884019917097105
$row[0,16] = 1.145619605068889
$row[0,17] = 10.31057644562
$row[0,18] = 0.01141722104999427
$row[0,19] = 0.01141722104999427
$ws.Range("A6:T6").Value = $row

# Row 7
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "FAPs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 0.3739756666666667
$row[0,7] = 1.121927
$row[0,8] = 0.06060032033836409
$row[0,9] = 0.06060032033836409
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 6.757657999999999
$row[0,13] = 20.272974
$row[0,14] = 0.4156086771445645
$row[0,15] = 0.4156086771445645
$row[0,16] = 2.527199655655333
$row[0,17] = 22.744796900898
$row[0,18] = 0.02518601897036434
$row[0,19] = 0.02518601897036434
$ws.Range("A7:T7").Value = $row

# Row 8
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "MuSCs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 0.3739756666666667
$row[0,7] = 1.121927
$row[0,8] = 0.06060032033836409
$row[0,9] = 0.06060032033836409
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.493414666666666
$row[0,13] = 10.480244
$row[0,14] = 0.214851572590793
$row[0,15] = 0.214851572590793
$row[0,16] = 1.306452078909778
$row[0,17] = 11.758068710188
$row[0,18] = 0.01302007412420334
$row[0,19] = 0.01302007412420334
$ws.Range("A8:T8").Value = $row

# Row 9
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "Resolving-Mac"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 0.3739756666666667
$row[0,7] = 1.121927
$row[0,8] = 0.06060032033836409
$row[0,9] = 0.06060032033836409
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.945239333333333
$row[0,13] = 8.835718
$row[0,14] = 0.1811377585549322
$row[0,15] = 0.1811377585549322
$row[0,16] = 1.101447843176222
$row[0,17] = 9.913030588586
$row[0,18] = 0.01097700619380214
$row[0,19] = 0.01097700619380214
$ws.Range("A9:T9").Value = $row

# Row 10
$row = New-Object 'object[,]' 1,20
$row[0,0] = "MuSCs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.624333
$row[0,7] = 13.872999
$row[0,8] = 0.7493430352008683
$row[0,9] = 0.7493430352008682
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.063353333333333
$row[0,13] = 9.190059999999999
$row[0,14] = 0.1884019917097105
$row[0,15] = 0.1884019917097105
$row[0,16] = 14.16596590999333
$row[0,17] = 127.49369318994
$row[0,18] = 0.1411777203056433
$row[0,19] = 0.1411777203056433
$ws.Range("A10:T10").Value = $row

# Row 11
$row = New-Object 'object[,]' 1,20
$row[0,0] = "MuSCs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.624333
$row[0,7] = 13.872999
$row[0,8] = 0.7493430352008683
$row[0,9] = 0.7493430352008682
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 6.757657999999999
$row[0,13] = 20.272974
$row[0,14] = 0.4156086771445645
$row[0,15] = 0.4156086771445645
$row[0,16] = 31.249660892114
$row[0,17] = 281.246948029026
$row[0,18] = 0.3114334675873257
$row[0,19] = 0.3114334675873257
$ws.Range("A11:T11").Value = $row

# Row 12
$row = New-Object 'object[,]' 1,20
$row[0,0] = "MuSCs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "MuSCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.624333
$row[0,7] = 13.872999
$row[0,8] = 0.7493430352008683
$row[0,9] = 0.7493430352008682
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.493414666666666
$row[0,13] = 10.480244
$row[0,14] = 0.214851572590793
$row[0,15] = 0.214851572590793
$row[0,16] = 16.15471272575067
$row[0,17] = 145.392414531756
$row[0,18] = 0.1609975295228645
$row[0,19] = 0.1609975295228645
$ws.Range("A12:T12").Value = $row

# Row 13
$row = New-Object 'object[,]' 1,20
$row[0,0] = "MuSCs"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "Resolving-Mac"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.624333
$row[0,7] = 13.872999
$row[0,8] = 0.7493430352008683
$row[0,9] = 0.7493430352008682
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.945239333333333
$row[0,13] = 8.835718
$row[0,14] = 0.1811377585549322
$row[0,15] = 0.1811377585549322
$row[0,16] = 13.61976744203133
$row[0,17] = 122.577906978282
$row[0,18] = 0.1357343177850349
$row[0,19] = 0.1357343177850349
$ws.Range("A13:T13").Value = $row

# Row 14
$row = New-Object 'object[,]' 1,20
$row[0,0] = "Resolving-Mac"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "ECs"
$row[0,4] = 1
$row[0,5] = 0.3333333333333333
$row[0,6] = 0.007355666666666667
$row[0,7] = 0.022067
$row[0,8] = 0.001191937861292829
$row[0,9] = 0.001191937861292829
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.063353333333333
$row[0,13] = 9.190059999999999
$row[0,14] = 0.1884019917097105
$row[0,15] = 0.1884019917097105
$row[0,16] = 0.02253300600222222
$row[0,17] = 0.20279705402
$row[0,18] = 0.0002245634670617817
$row[0,19] = 0.0002245634670617817
$ws.Range("A14:T14").Value = $row

# Row 15
$row = New-Object 'object[,]' 1,20
$row[0,0] = "Resolving-Mac"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "FAPs"
$row[0,4] = 1
$row[0,5] = 0.3333333333333333
$row[0,6] = 0.007355666666666667
$row[0,7] = 0.022067
$row[0,8] = 0.001191937861292829
$row[0,9] = 0.001191937861292829
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 6.757657999999999
$row[0,13] = 20.272974
$row[0,14] = 0.4156086771445645
$row[0,15] = 0.4156086771445645
$row[0,16] = 0.04970707969533333
$row[0,17] = 0.4473637172579999
$row[0,18] = 0.0004953797177704341
$row[0,19] = 0.0004953797177704341
$ws.Range("A15:T15").Value = $row

# Row 16
$row = New-Object 'object[,]' 1,20
$row[0,0] = "Resolving-Mac"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "MuSCs"
$row[0,4] = 1
$row[0,5] = 0.3333333333333333
$row[0,6] = 0.007355666666666667
$row[0,7] = 0.022067
$row[0,8] = 0.001191937861292829
$row[0,9] = 0.001191937861292829
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.493414666666666
$row[0,13] = 10.480244
$row[0,14] = 0.214851572590793
$row[0,15] = 0.214851572590793
$row[0,16] = 0.02569639381644444
$row[0,17] = 0.231267544348
$row[0,18] = 0.0002560897239292708
$row[0,19] = 0.0002560897239292708
$ws.Range("A16:T16").Value = $row

# Row 17
$row = New-Object 'object[,]' 1,20
$row[0,0] = "Resolving-Mac"
$row[0,1] = "Inhbb"
$row[0,2] = "Acvr1b"
$row[0,3] = "Resolving-Mac"
$row[0,4] = 1
$row[0,5] = 0.3333333333333333
$row[0,6] = 0.007355666666666667
$row[0,7] = 0.022067
$row[0,8] = 0.001191937861292829
$row[0,9] = 0.001191937861292829
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.945239333333333
$row[0,13] = 8.835718
$row[0,14] = 0.1811377585549322
$row[0,15] = 0.1811377585549322
$row[0,16] = 0.02166419878955556
$row[0,17] = 0.194977789106
$row[0,18] = 0.0002159049525313427
$row[0,19] = 0.0002159049525313427
$ws.Range("A17:T17").Value = $row
